# Actualizar precios con datos nuevos
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the timestamp on the last existing row (precision re-write, same instant)
$ws.Range("A8").Value = 45808.39130435186

# Append the new price entry as row 9
$ws.Range("A9").Value = 45809.39161978372
$ws.Range("A9").NumberFormat = $ws.Range("A8").NumberFormat

$ws.Range("B9").Value = "EVOWHEY PROTEIN"
$ws.Range("C9").Value = "2Kg"
$ws.Range("D9").Value = "33,90€"
